# Change tracing strategy and save wallet labels:
# append newly-traced wallet dates to the bottom of the Date column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @(
    "2024-08-28",
    "2024-09-21",
    "2024-09-09",
    "2024-09-06",
    "2024-10-04",
    "2024-10-03",
    "2024-10-01",
    "2024-10-05",
    "2024-08-31",
    "2024-09-17"
)

$startRow = 10
for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    # Force text interpretation (dates would otherwise auto-convert to
    # serial date numbers), then drop the formatting again so the cell
    # keeps the sheet's default (unstyled) appearance.
    $cell.NumberFormat = "@"
    $cell.Value = $newDates[$i]
    $cell.ClearFormats()
}
